# Updates IFRS financial figures for rows 2-6 (company_list sheet) to the revised
# (much smaller-scale) figures, and clears out rows 7-9 entirely except for the
# identifying columns A-C, matching the corrected dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "D2" = 20329
    "E2" = 74
    "F2" = 74
    "G2" = 786
    "H2" = 622
    "I2" = 607
    "J2" = 14
    "K2" = 16267
    "L2" = 7669
    "M2" = 8598
    "N2" = 8366
    "O2" = 232
    "P2" = 350
    "Q2" = 549
    "R2" = -1385
    "S2" = -120
    "T2" = 695
    "U2" = -146
    "V2" = 377
    "W2" = 0.37
    "X2" = 3.06
    "Y2" = 7.52
    "Z2" = 3.94
    "AA2" = 89.2
    "AB2" = 2429.42
    "AC2" = 12145
    "AD2" = 11.77
    "AE2" = 215327
    "AF2" = 0.66
    "AG2" = 1500
    "AH2" = 1.05
    "AI2" = 9.6
    "AJ2" = 5000000
    "D3" = 15403
    "E3" = -25
    "F3" = -25
    "G3" = 554
    "H3" = 448
    "I3" = 445
    "J3" = 3
    "K3" = 14366
    "L3" = 5625
    "M3" = 8741
    "N3" = 8694
    "O3" = 47
    "P3" = 350
    "Q3" = 103
    "R3" = 431
    "S3" = -97
    "T3" = 448
    "U3" = -345
    "V3" = 125
    "W3" = -0.16
    "X3" = 2.91
    "Y3" = 5.22
    "Z3" = 2.92
    "AA3" = 64.34999999999999
    "AB3" = 2523.28
    "AC3" = 8900
    "AD3" = 12.75
    "AE3" = 223759
    "AF3" = 0.51
    "AG3" = 1750
    "AH3" = 1.54
    "AI3" = 15.28
    "AJ3" = 5000000
    "D4" = 12720
    "E4" = -67
    "F4" = -67
    "G4" = 330
    "H4" = 235
    "I4" = 232
    "J4" = 3
    "K4" = 14666
    "L4" = 5761
    "M4" = 8905
    "N4" = 8863
    "O4" = 42
    "P4" = 350
    "Q4" = 768
    "R4" = -986
    "S4" = -59
    "T4" = 399
    "U4" = 370
    "V4" = 138
    "W4" = -0.53
    "X4" = 1.85
    "Y4" = 2.64
    "Z4" = 1.62
    "AA4" = 64.69
    "AB4" = 2571.92
    "AC4" = 4640
    "AD4" = 17.72
    "AE4" = 228122
    "AF4" = 0.36
    "AG4" = 1750
    "AH4" = 2.13
    "AI4" = 29.31
    "AJ4" = 5000000
    "D5" = 13504
    "E5" = 101
    "F5" = 101
    "G5" = 477
    "H5" = 366
    "I5" = 363
    "J5" = 2
    "K5" = 15459
    "L5" = 6293
    "M5" = 9165
    "N5" = 9124
    "O5" = 41
    "P5" = 350
    "Q5" = 448
    "R5" = -273
    "S5" = -77
    "T5" = 272
    "U5" = 176
    "V5" = 132
    "W5" = 0.75
    "X5" = 2.71
    "Y5" = 4.04
    "Z5" = 2.43
    "AA5" = 68.66
    "AB5" = 2645.84
    "AC5" = 7265
    "AD5" = 12.98
    "AE5" = 234861
    "AF5" = 0.4
    "AG5" = 1750
    "AH5" = 1.86
    "AI5" = 18.72
    "AJ5" = 5000000
    "D6" = 13890
    "E6" = 43
    "F6" = 43
    "G6" = 489
    "H6" = 418
    "I6" = 414
    "K6" = 13686
    "L6" = 4211
    "M6" = 9475
    "N6" = 9436
    "P6" = 350
    "Q6" = -1529
    "R6" = 929
    "S6" = -75
    "T6" = 281
    "U6" = -1810
    "V6" = 135
    "W6" = 0.31
    "X6" = 3.01
    "Y6" = 4.46
    "Z6" = 2.87
    "AA6" = 44.45
    "AB6" = 2736.72
    "AC6" = 8276
    "AD6" = 10.32
    "AE6" = 242938
    "AF6" = 0.35
    "AG6" = 1750
    "AH6" = 2.05
    "AI6" = 16.43
    "AJ6" = 5000000
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}

# Rows 7-9 (companies #40, #41, #42) no longer carry any of the financial-metric
# columns in the corrected data set - only the leading identifier columns (A:C) remain.
$ws.Range("D7:AJ9").ClearContents()
